$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "XGBOOST" algorithm row (row 8) to the metrics table.
$ws.Range("A8").Value = "XGBOOST"
$ws.Range("B8").Value = 0.9037
$ws.Range("C8").Value = 47.14
$ws.Range("D8").Value = 4133.21
$ws.Range("E8").Value = 64.29

# Move the active selection (matches the saved cursor position in the file).
$ws.Range("G7").Select() | Out-Null
